$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Comparaciones_Significativas (col B) and Proporcion_Sig (col C)
# Row 4: AV-MCPS
$ws.Range("B4").Value = "3/10"
$ws.Range("C4").Value = 76.8

# Row 5: DeepAR
$ws.Range("B5").Value = "2/10"
$ws.Range("C5").Value = 51.2

# Row 6: Sieve Bootstrap
$ws.Range("B6").Value = "1/10"
$ws.Range("C6").Value = 25.6
$ws.Range("E6").Value = 0.8531623364799357

# Row 7: AREPD
$ws.Range("B7").Value = "0/10"
$ws.Range("C7").Value = 0

# Row 8: Block Bootstrapping
$ws.Range("B8").Value = "0/10"
$ws.Range("C8").Value = 0
$ws.Range("E8").Value = 0.7552480947305171

# Row 9: LSPMW
$ws.Range("B9").Value = "0/10"
$ws.Range("C9").Value = 0

# Row 10: LSPM
$ws.Range("B10").Value = "0/10"
$ws.Range("C10").Value = 0
$ws.Range("E10").Value = 0.5327762422047049
